# Scheduled-runner price refresh: update market/profit figures (columns H:N)
# on a handful of Leve rows across several crafting-job sheets, and blank out
# the stale H:N figures for a block of fully-delisted rows on CRP.
#
# Columns: H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#          K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ---- ALC ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H113").Value = 2337.5
$ws.Range("I113").Value = 2202.1428
$ws.Range("J113").Value = 2653.3333
$ws.Range("K113").Value = 2202.1428
$ws.Range("L113").Value = 2653.3333
$ws.Range("M113").Value = 1051.8572
$ws.Range("N113").Value = -9161.3333

$ws.Range("H132").Value = 2679.3015
$ws.Range("I132").Value = 2323.9783
$ws.Range("J132").Value = 3640.7646
$ws.Range("K132").Value = 6971.9349
$ws.Range("L132").Value = 10922.2938
$ws.Range("M132").Value = -4441.9349
$ws.Range("N132").Value = -15982.2938

$ws.Range("H137").Value = 5315.793
$ws.Range("I137").Value = 1394
$ws.Range("J137").Value = 6809.8096
$ws.Range("K137").Value = 4182
$ws.Range("L137").Value = 20429.4288
$ws.Range("M137").Value = -1632
$ws.Range("N137").Value = -25529.4288

$ws.Range("H138").Value = 2540.258
$ws.Range("I138").Value = 1177.8182
$ws.Range("J138").Value = 5870.6665
$ws.Range("K138").Value = 3533.4546
$ws.Range("L138").Value = 17611.9995
$ws.Range("M138").Value = 1606.5454
$ws.Range("N138").Value = -27891.9995

# ---- ARM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H132").Value = 13359.632
$ws.Range("I132").Value = 11721.182
$ws.Range("J132").Value = 15612.5
$ws.Range("K132").Value = 35163.546
$ws.Range("L132").Value = 46837.5
$ws.Range("M132").Value = -32633.546
$ws.Range("N132").Value = -51897.5

# ---- BSM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 2276.2942
$ws.Range("I134").Value = 2011.3549
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 6034.0647
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -3499.0647
$ws.Range("N134").Value = -20112

# ---- CRP --------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 6464.5454
$ws.Range("I31").Value = 2080.7273
$ws.Range("K31").Value = 2080.7273
$ws.Range("M31").Value = -1785.7273

$ws.Range("H34").Value = 6464.5454
$ws.Range("I34").Value = 2080.7273
$ws.Range("K34").Value = 2080.7273
$ws.Range("M34").Value = -1878.7273

$ws.Range("H93").Value = 8802.333000000001
$ws.Range("I93").Value = 8802.333000000001
$ws.Range("K93").Value = 8802.333000000001
$ws.Range("M93").Value = -6930.333000000001

$ws.Range("H103").Value = 6342
$ws.Range("I103").Value = 6342
$ws.Range("K103").Value = 6342
$ws.Range("M103").Value = -5170

# Rows 129-135 and 137-141: these Leve items no longer have market data, so
# the whole H:N block is cleared out (row 136 is left untouched).
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("N129").ClearContents()

$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()

$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("N131").ClearContents()

$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()

$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()

$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("N138").ClearContents()

$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("N139").ClearContents()

$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("N140").ClearContents()

$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

# ---- CUL --------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 482.36
$ws.Range("I5").Value = 433
$ws.Range("J5").Value = 1050
$ws.Range("K5").Value = 1299
$ws.Range("L5").Value = 3150
$ws.Range("M5").Value = -1187
$ws.Range("N5").Value = -3374

$ws.Range("H122").Value = 718.7308
$ws.Range("I122").Value = 329
$ws.Range("J122").Value = 1052.7858
$ws.Range("K122").Value = 2961
$ws.Range("L122").Value = 9475.072200000001
$ws.Range("M122").Value = -511
$ws.Range("N122").Value = -14375.0722

$ws.Range("H132").Value = 2440082.8
$ws.Range("I132").Value = 4762765
$ws.Range("J132").Value = 1266.25
$ws.Range("K132").Value = 42864885
$ws.Range("L132").Value = 11396.25
$ws.Range("M132").Value = -42862355
$ws.Range("N132").Value = -16456.25

$ws.Range("H135").Value = 482.36
$ws.Range("I135").Value = 433
$ws.Range("J135").Value = 1050
$ws.Range("K135").Value = 3897
$ws.Range("L135").Value = 9450
$ws.Range("M135").Value = -1362
$ws.Range("N135").Value = -14520

# ---- GSM --------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H132").Value = 6027.8857
$ws.Range("I132").Value = 6877.769
$ws.Range("J132").Value = 3572.6667
$ws.Range("K132").Value = 20633.307
$ws.Range("L132").Value = 10718.0001
$ws.Range("M132").Value = -18103.307
$ws.Range("N132").Value = -15778.0001

# ---- LTW --------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1522.75
$ws.Range("I22").Value = 1153.75
$ws.Range("J22").Value = 1707.25
$ws.Range("K22").Value = 1153.75
$ws.Range("L22").Value = 1707.25
$ws.Range("M22").Value = -858.75
$ws.Range("N22").Value = -2297.25

$ws.Range("H27").Value = 1522.75
$ws.Range("I27").Value = 1153.75
$ws.Range("J27").Value = 1707.25
$ws.Range("K27").Value = 1153.75
$ws.Range("L27").Value = 1707.25
$ws.Range("M27").Value = -1046.75
$ws.Range("N27").Value = -1921.25

$ws.Range("H46").Value = 1456.2307
$ws.Range("I46").Value = 2857.75
$ws.Range("J46").Value = 833.3333
$ws.Range("K46").Value = 2857.75
$ws.Range("L46").Value = 833.3333
$ws.Range("M46").Value = -2669.75
$ws.Range("N46").Value = -1209.3333

# ---- WVR --------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

$ws.Range("H132").Value = 5053.2
$ws.Range("I132").Value = 5612.448
$ws.Range("J132").Value = 3578.818
$ws.Range("K132").Value = 16837.344
$ws.Range("L132").Value = 10736.454
$ws.Range("M132").Value = -14307.344
$ws.Range("N132").Value = -15796.454
